# Insert a new weekly record at row 31 (Terminal Hortofrutícola Agro Chillán - Mango),
# shifting the existing rows 31-75 down to 32-76.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(31).Insert()

$ws.Range("A31").Value = 7
$ws.Range("B31").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C31").Value = "Ñuble"
$ws.Range("D31").Value = 44650
$ws.Range("E31").Value = 16
$ws.Range("F31").Value = "Fruta"
$ws.Range("G31").Value = 100108
$ws.Range("H31").Value = "Tropicales y subtropicales"
$ws.Range("I31").Value = 100108002
$ws.Range("J31").Value = "Mango"
$ws.Range("K31").Value = "Sin especificar"
$ws.Range("L31").Value = "Primera"
$ws.Range("M31").Value = 120
$ws.Range("N31").Value = 7500
$ws.Range("O31").Value = 8000
$ws.Range("P31").Value = 7750
$ws.Range("Q31").Value = "$/bandeja 4 kilos"
$ws.Range("R31").Value = "Perú"
$ws.Range("S31").Value = 1938
$ws.Range("T31").Value = 4
